$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 0.1192916666666667
$ws.Range("H2").Value = 0.357875
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 0.9623583333333334
$ws.Range("N2").Value = 2.887075
$ws.Range("O2").Value = 0.05805926999654511
$ws.Range("P2").Value = 0.05805926999654509
$ws.Range("Q2").Value = 0.1148013295138889
$ws.Range("R2").Value = 1.033211965625
$ws.Range("S2").Value = 0.05805926999654511
$ws.Range("T2").Value = 0.05805926999654509

# Row 3
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 0.1192916666666667
$ws.Range("H3").Value = 0.357875
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 12.12890133333333
$ws.Range("N3").Value = 36.386704
$ws.Range("O3").Value = 0.731739034081334
$ws.Range("P3").Value = 0.7317390340813339
$ws.Range("Q3").Value = 1.446876854888889
$ws.Range("R3").Value = 13.021891694
$ws.Range("S3").Value = 0.731739034081334
$ws.Range("T3").Value = 0.7317390340813339

# Row 4
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 0.1192916666666667
$ws.Range("H4").Value = 0.357875
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 3.484187
$ws.Range("N4").Value = 10.452561
$ws.Range("O4").Value = 0.210201695922121
$ws.Range("P4").Value = 0.2102016959221209
$ws.Range("Q4").Value = 0.4156344742083333
$ws.Range("R4").Value = 3.740710267875
$ws.Range("S4").Value = 0.210201695922121
$ws.Range("T4").Value = 0.2102016959221209
